{"js": "// Remove the trailing \"Ver no Jupiter ...\" / site-footer (\"\u00a9 2020 ...\")\n// paragraphs (and the blank paragraph that separates them from the\n// bibliography) that were scraped onto the end of the document, while\n// leaving the final blank / page-break paragraphs untouched.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Anchor on the last bibliography entry so the edit is resilient to the\n// exact paragraph count/position rather than relying on hard-coded indices.\nconst anchorMarker = \"Thomson Pioneira\";\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorMarker) !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the bibliography anchor paragraph.\");\n}\n\n// The three paragraphs that follow the anchor are the blank separator,\n// the \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" footer line - all of\n// which must go. Delete from the highest index down so earlier indices\n// stay valid.\nconst toRemove = [];\nfor (let i = anchorIndex + 1; i <= anchorIndex + 3 && i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (\n    text.trim() === \"\" ||\n    text.indexOf(\"Ver no Jupiter\") !== -1 ||\n    text.indexOf(\"Contact: luizeleno@usp.br\") !== -1\n  ) {\n    toRemove.push(i);\n  }\n}\n\nfor (let i = toRemove.length - 1; i >= 0; i--) {\n  paragraphs.items[toRemove[i]].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / site-footer (\"\u00a9 2020 ...\")\n# paragraphs (and the blank paragraph that separates them from the\n# bibliography) that were scraped onto the end of the document, while\n# leaving the final blank / page-break paragraphs untouched.\n$d = $word.ActiveDocument\n\n# Anchor on the last bibliography entry so the edit is resilient to the\n# exact paragraph count/position rather than relying on hard-coded indices.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Thomson Pioneira*\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the bibliography anchor paragraph.\"\n}\n\n# The three paragraphs that follow the anchor are the blank separator,\n# the \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" footer line - all of\n# which must go. Collect them first, then delete starting from the\n# highest index so earlier indices stay valid.\n$toRemove = @()\n$upper = [Math]::Min($anchorIndex + 3, $d.Paragraphs.Count)\nfor ($i = $anchorIndex + 1; $i -le $upper; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($text -eq \"\" -or $text -like \"*Ver no Jupiter*\" -or $text -like \"*Contact: luizeleno@usp.br*\") {\n        $toRemove += $i\n    }\n}\n\nfor ($j = $toRemove.Count - 1; $j -ge 0; $j--) {\n    $d.Paragraphs.Item($toRemove[$j]).Range.Delete()\n}\n"}
